# The source data contained embedded line breaks inside several "Address" /
# "Activities" / "Details" cells (e.g. "عنوان`nالوليد ب..."). This collapses
# each embedded newline into a single space so the text reads on one line
# (e.g. "عنوان الوليد ب...").

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Sheet1-2024")
$ws.Range('K2').Value  = 'عنوان الوليد ب...'
$ws.Range('K3').Value  = 'عنوان 3423 سيل...'
$ws.Range('K4').Value  = 'عنوان الدور ال...'
$ws.Range('M4').Value  = 'الدور الأول ٧٧...'
$ws.Range('K5').Value  = 'عنوان الرياض ش...'
$ws.Range('K6').Value  = 'عنوان 4556 طري...'
$ws.Range('K7').Value  = 'عنوان منطقة جا...'
$ws.Range('K9').Value  = 'عنوان 7827 عرف...'
$ws.Range('K10').Value = 'عنوان طريق الا...'
$ws.Range('K11').Value = 'عنوان الدمام /...'
$ws.Range('K13').Value = 'عنوان مبنى رقم...'
$ws.Range('K18').Value = 'عنوان الغيثي ب...'

# "تشييد المباني`n..." -> "تشييد المباني ..." repeats across every row of
# the Activities column (N) that shared that value.
foreach ($cell in @('N2','N4','N5','N6','N7','N9','N12','N13','N15','N20','N23')) {
    $ws.Range($cell).Value = 'تشييد المباني ...'
}

$ws = $wb.Worksheets.Item("Sheet2-2022")
$ws.Range('K2').Value  = '   شركة المنار...'
$ws.Range('K3').Value  = '   تأسست الشرك...'
$ws.Range('K5').Value  = '   مقابل مكتبة...'
$ws.Range('K6').Value  = '   نحن متخصصون...'
$ws.Range('K15').Value = '   طريق الملك ...'
$ws.Range('K19').Value = '   شارع الأمير...'
$ws.Range('K21').Value = '   بجانب البنك...'
$ws.Range('K24').Value = '   البكيرية-طر...'
